# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 16 / Row 17: swap "Periodo Mora" (E) and "Valor Mora" (F) values between
# the two detail rows, and update "Salario Basico" (G) for both rows.
$ws.Range("E16").Value = "2309"
$ws.Range("F16").Value = 58667
$ws.Range("G16").Value = 2000000

$ws.Range("E17").Value = "2310"
$ws.Range("F17").Value = 80000
$ws.Range("G17").Value = 2000000
